$d = $word.ActiveDocument

# ------------------------------------------------------------------
# 1) Weekly date bump: "13 October 2014" -> "22 October 2014"
#    (cached SAVEDATE field result text)
# ------------------------------------------------------------------
$d.Content.Find.Execute("13 October 2014", $true, $false, $false, $false, `
                         $false, $true, 1, $false, "22 October 2014", 2) | Out-Null

# ------------------------------------------------------------------
# 2) Tag the three screenshot runs (InlineShapes) with en-GB language
#    so their <w:rPr> gains <w:lang w:val="en-GB"/> next to <w:noProof/>
# ------------------------------------------------------------------
for ($i = 1; $i -le $d.InlineShapes.Count; $i++) {
    $shp = $d.InlineShapes.Item($i)
    $shp.Range.LanguageID = "en-GB"
}

# ------------------------------------------------------------------
# 3) Append the new "Log and Error Handling" section after the
#    existing "From Source tree ..." paragraph, moving the _GoBack
#    bookmark to the end of the freshly-added content.
# ------------------------------------------------------------------
$goBack = $d.Bookmarks.Item("_GoBack")
$goBack.Delete()

$lastPara = $d.Paragraphs.Item($d.Paragraphs.Count)
$paraEnd = $lastPara.Range.End
$insertPoint = $d.Range($paraEnd, $paraEnd)

$newContent = '<w:p><w:pPr><w:pStyle w:val="Heading1"/></w:pPr>' + `
    '<w:r><w:t xml:space="preserve">Log and </w:t></w:r>' + `
    '<w:r><w:t>Error Handling</w:t></w:r></w:p>' + `
    '<w:p><w:r><w:t xml:space="preserve">Error handing works close with Logger and </w:t></w:r>' + `
    '<w:proofErr w:type="spellStart"/><w:r><w:t>LogCode</w:t></w:r><w:proofErr w:type="spellEnd"/>' + `
    '<w:r><w:t>. Depend on log code profile to decide which to log</w:t></w:r></w:p>' + `
    '<w:p><w:r><w:t>Log Code Category</w:t></w:r></w:p>' + `
    '<w:p><w:r><w:t>100-199</w:t></w:r>' + `
    '<w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/></w:p>'

$xml = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?>' + `
    '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' + `
    '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' + `
    '<pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body>' + `
    $newContent + `
    '</w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'

$insertPoint.InsertXML($xml) | Out-Null

Write-Host "edit complete"
